$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: id=3, nombre=(empty), edad=14
$ws.Range("A4").Value = 3
$ws.Range("C4").Value = 14

# Row 5: id=4, nombre=alberto, edad=18
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "alberto"
$ws.Range("C5").Value = 18

[void]$ws.Range("C5").Select()
